# Generate Report for Handoff
# Adds two new localized-file rows (the .png assets) to the Overview sheet
# and to each language sheet (zh-cn / de-de), and refreshes the existing
# "bb741da7-...md" handoff row to point at the new "0098be31-...md" handoff.

$wb = $excel.ActiveWorkbook

$guidOld  = "bb741da7-1861-47e3-aa4e-9d3bd5b6c2ee"
$guidNew  = "0098be31-459f-4b6e-9207-3c6d49873c86"
$hashNew  = "3085e53c027ba54fb3c094d237a8d37014fd6b04"

$mdName      = "$guidNew.md"
$zhXlfName   = "$guidNew.$hashNew.zh-cn.xlf"
$deXlfName   = "$guidNew.$hashNew.de-de.xlf"

$png1Name    = "495d83d2-554f-4e75-bf1d-1879f1190e36.png"
$png2Name    = "516faded-d99b-45c3-ad0c-b3d423d6807f.png"
$png1Target  = "2ca720aa2f512d03097ec5e37bfa81da43ca7ead.png"
$png2Target  = "44472dc3c011469ee3f4cf8e09eeab5e6297abfe.png"

$statusReady   = "Ready for handoff"
$statusInclude = "Include"
$statusDep     = "IsDependency"
$zeroDate      = "0001-01-01 00:00:00"
$dependencyFrom = "e2e\$mdName"

$overviewDate  = "2016-47-19 10:47:43"
$zhHandoffDate = "2016-03-19 10:47:40"
$deHandoffDate = "2016-03-19 10:47:43"

$mdUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/803170674b95a6668190645c8d413c30d264d69b/e2e/$mdName"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7a1fdf5faa66375768d4f78c33b969b2b3ae6eff/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlfName"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fc0d4b1b690858d76c2efb6b84a77f1adf77b850/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlfName"

$png1Url      = "https://github.com/OpenLocalizationTest/oltest/blob/803170674b95a6668190645c8d413c30d264d69b/e2e/$png1Name"
$png2Url      = "https://github.com/OpenLocalizationTest/oltest/blob/803170674b95a6668190645c8d413c30d264d69b/e2e/$png2Name"
$png1ExtUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/803170674b95a6668190645c8d413c30d264d69b/e2e/$png1Name"
$png2ExtUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/803170674b95a6668190645c8d413c30d264d69b/e2e/$png2Name"
$png1HtZhUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7a1fdf5faa66375768d4f78c33b969b2b3ae6eff/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$png1Target"
$png2HtZhUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7a1fdf5faa66375768d4f78c33b969b2b3ae6eff/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$png2Target"
$png1HtDeUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fc0d4b1b690858d76c2efb6b84a77f1adf77b850/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$png1Target"
$png2HtDeUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fc0d4b1b690858d76c2efb6b84a77f1adf77b850/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$png2Target"

function Set-Link($ws, $cellRef, $text, $url) {
    $rng = $ws.Range($cellRef)
    if ($rng.Hyperlinks.Count -gt 0) {
        $rng.Hyperlinks.Delete()
    }
    $ws.Hyperlinks.Add($rng, $url, "", "", $text) | Out-Null
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

Set-Link $ws1 "A2" $mdName $mdUrl

$ws1.Range("B2").Value = $statusReady
$ws1.Range("C2").Value = $statusReady
$ws1.Range("D2").Value = $overviewDate

Set-Link $ws1 "A3" $png1Name $png1Url
$ws1.Range("B3").Value = $statusReady
$ws1.Range("C3").Value = $statusReady
$ws1.Range("D3").Value = $overviewDate

Set-Link $ws1 "A4" $png2Name $png2Url
$ws1.Range("B4").Value = $statusReady
$ws1.Range("C4").Value = $statusReady
$ws1.Range("D4").Value = $overviewDate

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

Set-Link $ws2 "A2" $mdName $mdUrl
Set-Link $ws2 "B2" ".md" $mdUrl
$ws2.Range("C2").Value = $statusReady
Set-Link $ws2 "D2" $zhXlfName $zhXlfUrl
$ws2.Range("E2").Value = $zhHandoffDate
$ws2.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H2").Value = $zeroDate
$ws2.Range("I2").Value = $statusInclude

Set-Link $ws2 "A3" $png1Name $png1Url
Set-Link $ws2 "B3" ".png" $png1ExtUrl
$ws2.Range("C3").Value = $statusReady
Set-Link $ws2 "D3" $png1Target $png1HtZhUrl
$ws2.Range("E3").Value = $zhHandoffDate
$ws2.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H3").Value = $zeroDate
$ws2.Range("I3").Value = $statusDep
$ws2.Range("J3").Value = $dependencyFrom

Set-Link $ws2 "A4" $png2Name $png2Url
Set-Link $ws2 "B4" ".png" $png2ExtUrl
$ws2.Range("C4").Value = $statusReady
Set-Link $ws2 "D4" $png2Target $png2HtZhUrl
$ws2.Range("E4").Value = $zhHandoffDate
$ws2.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H4").Value = $zeroDate
$ws2.Range("I4").Value = $statusDep
$ws2.Range("J4").Value = $dependencyFrom

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

Set-Link $ws3 "A2" $mdName $mdUrl
Set-Link $ws3 "B2" ".md" $mdUrl
$ws3.Range("C2").Value = $statusReady
Set-Link $ws3 "D2" $deXlfName $deXlfUrl
$ws3.Range("E2").Value = $deHandoffDate
$ws3.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H2").Value = $zeroDate
$ws3.Range("I2").Value = $statusInclude

Set-Link $ws3 "A3" $png1Name $png1Url
Set-Link $ws3 "B3" ".png" $png1ExtUrl
$ws3.Range("C3").Value = $statusReady
Set-Link $ws3 "D3" $png1Target $png1HtDeUrl
$ws3.Range("E3").Value = $deHandoffDate
$ws3.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H3").Value = $zeroDate
$ws3.Range("I3").Value = $statusDep
$ws3.Range("J3").Value = $dependencyFrom

Set-Link $ws3 "A4" $png2Name $png2Url
Set-Link $ws3 "B4" ".png" $png2ExtUrl
$ws3.Range("C4").Value = $statusReady
Set-Link $ws3 "D4" $png2Target $png2HtDeUrl
$ws3.Range("E4").Value = $deHandoffDate
$ws3.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H4").Value = $zeroDate
$ws3.Range("I4").Value = $statusDep
$ws3.Range("J4").Value = $dependencyFrom

Write-Output "Report generated."
